# Add the new "Anna Zandonati" team as row 64 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Value = "Anna Zandonati"
$ws.Range("B64").Value = "Stefano Tita | Clitoriders"
$ws.Range("C64").Value = "Mattia Festi | Shark Attack"
$ws.Range("D64").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E64").Value = "Federico  Manica | iMontagna"
$ws.Range("F64").Value = "Alessandro  Tengattini | Herta Vernello"
